$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the schedule day/slot cells (columns F, G, H, I) and the
#     credit count in D4, matching the new .csv source data ---

# Row 1 (CS02 / KTLT)
$ws.Cells.Item(1,6).Value = "Wed"
$ws.Cells.Item(1,8).Value = "Sat"
$ws.Cells.Item(1,9).Value = "S4"

# Row 2 (PHY01 / VLDC)
$ws.Cells.Item(2,6).Value = "Mon"
$ws.Cells.Item(2,7).Value = "S1"
$ws.Cells.Item(2,8).Value = "Mon"

# Row 3 (MAT02 / VTP1B)
$ws.Cells.Item(3,6).Value = "Tue"
$ws.Cells.Item(3,7).Value = "S2"
$ws.Cells.Item(3,8).Value = "Fri"
$ws.Cells.Item(3,9).Value = "S2"

# Row 4 (MLN01 / MLN)
$ws.Cells.Item(4,4).Value = 3
$ws.Cells.Item(4,6).Value = "Mon"
$ws.Cells.Item(4,7).Value = "S3"
$ws.Cells.Item(4,8).Value = "Mon"

# --- Resize columns A and C to fit their (now shorter) contents ---
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(3).AutoFit()

# --- Move the active selection to G4 ---
$ws.Range("G4").Select()
